$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Dispatchable_2023"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Nuclear_2023"
